$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ATATATATAT', 'WSWSWSWSWS', 0, 0, 0, $null, $null),
    @('ATAAATATAT', 'WSWSWSWSWS', 1, 0, 1, "'3", $null),
    @('ATATTAATAT', 'WSWSWSWSWS', 0, 1, 1, $null, "'5"),
    @('ATATAAATAT', 'WSWSWSWSWS', 1, 0, 1, "'5", $null),
    @('ATATATTAAT', 'WSWSWSWSWS', 0, 1, 1, $null, "'7"),
    @('ATTAATATAT', 'WSWSWSWSWS', 0, 1, 1, $null, "'3"),
    @('AAATATATAT', 'WSWSWSWSWS', 1, 0, 1, "'1", $null),
    @('TAATATATAT', 'WSWSWSWSWS', 0, 1, 1, $null, "'1"),
    @('ATATATAAAT', 'WSWSWSWSWS', 1, 0, 1, "'7", $null),
    @('ATATAAAAAT', 'WSWSWSWSWS', 2, 0, 2, "'5, 7", $null),
    @('AAATAAATAT', 'WSWSWSWSWS', 2, 0, 2, "'1, 5", $null),
    @('ATTAAAATAT', 'WSWSWSWSWS', 1, 1, 2, "'5", "'3"),
    @('ATAAATAAAT', 'WSWSWSWSWS', 2, 0, 2, "'3, 7", $null),
    @('TAAAATATAT', 'WSWSWSWSWS', 1, 1, 2, "'3", "'1"),
    @('AAATATTAAT', 'WSWSWSWSWS', 1, 1, 2, "'1", "'7"),
    @('ATAAATTAAT', 'WSWSWSWSWS', 1, 1, 2, "'3", "'7"),
    @('ATAATAATAT', 'WSWSWSWSWS', 1, 1, 2, "'3", "'5"),
    @('ATATAATAAT', 'WSWSWSWSWS', 1, 1, 2, "'5", "'7"),
    @('AATAATATAT', 'WSWSWSWSWS', 1, 1, 2, "'1", "'3"),
    @('TAATATAAAT', 'WSWSWSWSWS', 1, 1, 2, "'7", "'1"),
    @('ATAAAAATAT', 'WSWSWSWSWS', 2, 0, 2, "'3, 5", $null),
    @('AAATTAATAT', 'WSWSWSWSWS', 1, 1, 2, "'1", "'5"),
    @('AAATATAAAT', 'WSWSWSWSWS', 2, 0, 2, "'1, 7", $null),
    @('AAAAATATAT', 'WSWSWSWSWS', 2, 0, 2, "'1, 3", $null),
    @('TAATAAATAT', 'WSWSWSWSWS', 1, 1, 2, "'5", "'1"),
    @('ATAAAATAAT', 'WSWSWSWSWS', 2, 1, 3, "'3, 5", "'7"),
    @('AAAAATTAAT', 'WSWSWSWSWS', 2, 1, 3, "'1, 3", "'7"),
    @('AAAATAATAT', 'WSWSWSWSWS', 2, 1, 3, "'1, 3", "'5"),
    @('TAAAATAAAT', 'WSWSWSWSWS', 2, 1, 3, "'3, 7", "'1"),
    @('AAATAAAAAT', 'WSWSWSWSWS', 3, 0, 3, "'1, 5, 7", $null),
    @('ATAAAAAAAT', 'WSWSWSWSWS', 3, 0, 3, "'3, 5, 7", $null),
    @('AAATAATAAT', 'WSWSWSWSWS', 2, 1, 3, "'1, 5", "'7"),
    @('TAATAAAAAT', 'WSWSWSWSWS', 2, 1, 3, "'5, 7", "'1"),
    @('AAAAATAAAT', 'WSWSWSWSWS', 3, 0, 3, "'1, 3, 7", $null),
    @('AATAAAATAT', 'WSWSWSWSWS', 2, 1, 3, "'1, 5", "'3"),
    @('TAAAAAATAT', 'WSWSWSWSWS', 2, 1, 3, "'3, 5", "'1"),
    @('AAAAAAATAT', 'WSWSWSWSWS', 3, 0, 3, "'1, 3, 5", $null),
    @('AAAAAATAAT', 'WSWSWSWSWS', 3, 1, 4, "'1, 3, 5", "'7"),
    @('AAAAAAAAAT', 'WSWSWSWSWS', 4, 0, 4, "'1, 3, 5, 7", $null),
    @('TAAAAAAAAT', 'WSWSWSWSWS', 3, 1, 4, "'3, 5, 7", "'1")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}